# Apply the 2025-11-17 09:19 data refresh to market_health_data.xlsx
#
# Summary of changes:
#  1. Metadata!A2        - bump "Last Updated" timestamp by one minute.
#  2. Top Gainers        - a new top gainer (GROWW) enters the list at row 24;
#                          every following row shifts down by one and the
#                          previous last row (AWFIS) drops off the bottom.
#  3. Industry Analysis!C36 - "Latest" figure for the
#                          "automobiles - lcvs / hcvs" industry is corrected.
#  4. Stock List         - two new stocks (GROWW, TMCV) enter at rows 2-3;
#                          every following row shifts down by two and the
#                          previous last two rows (BLUSPRING, DIGITIDE) drop
#                          off the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata: "Last Updated" timestamp
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "17 Nov 2025, 09:18 AM"

# ---------------------------------------------------------------------
# 2. Top Gainers: insert GROWW as the new row 24, shift rows 24-75 down
#    into 25-76, and drop what used to be row 76 (AWFIS).
# ---------------------------------------------------------------------
$wsTG = $wb.Worksheets.Item("Top Gainers")
$wsTG.Range("A76").EntireRow.Delete()
$wsTG.Range("A24").EntireRow.Insert()
$wsTG.Range("A24:E24").ClearFormats()

$wsTG.Range("A24").Value = $wsTG.Range("A25").Value2
$wsTG.Range("B24").Value = "GROWW"
$wsTG.Range("C24").Value = 3.7837
$wsTG.Range("D24").Value = "N/A"
$wsTG.Range("E24").Value = "N/A"

# ---------------------------------------------------------------------
# 3. Industry Analysis: fix the "Latest" value for row 36
#    (automobiles - lcvs / hcvs) to match the Weekly column.
# ---------------------------------------------------------------------
$wsIA = $wb.Worksheets.Item("Industry Analysis")
$wsIA.Range("C36").Value = 0.4762

# ---------------------------------------------------------------------
# 4. Stock List: insert GROWW and TMCV as the new rows 2-3, shift rows
#    2-74 down into 4-76, and drop what used to be rows 75-76
#    (BLUSPRING, DIGITIDE).
# ---------------------------------------------------------------------
$wsSL = $wb.Worksheets.Item("Stock List")
$wsSL.Range("A75:A76").EntireRow.Delete()
$wsSL.Range("A2:A3").EntireRow.Insert()
$wsSL.Range("A2:H3").ClearFormats()

$wsSL.Range("A2").Value = $wsSL.Range("A4").Value2
$wsSL.Range("B2").Value = "GROWW"
$wsSL.Range("C2").Value = "GROWW"
$wsSL.Range("D2").Value = 154.15
$wsSL.Range("E2").Value = 3.7837
$wsSL.Range("F2").Value = "N/A"
$wsSL.Range("G2").Value = "N/A"
$wsSL.Range("H2").Value = 91696.4308

$wsSL.Range("A3").Value = $wsSL.Range("A4").Value2
$wsSL.Range("B3").Value = "TMCV"
$wsSL.Range("C3").Value = "TMCV"
$wsSL.Range("D3").Value = 321.8
$wsSL.Range("E3").Value = 1.3224
$wsSL.Range("F3").Value = "N/A"
$wsSL.Range("G3").Value = "N/A"
$wsSL.Range("H3").Value = 116950.8444
